# Add a new "2022" column (Q) to the table, mirroring the formatting of
# the existing "2021" column (P). Rows 10-25 repeat the "…" placeholder
# string already used in column P; rows 3, 6 and 9 stay blank (only
# inheriting the border/format style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (and, for the placeholder rows, the value) from column P
# into the new column Q for every row of the table.
for ($r = 3; $r -le 25; $r++) {
    $src = $ws.Cells.Item($r, 16)   # column P
    $dst = $ws.Cells.Item($r, 17)   # column Q
    $src.Copy($dst)
}

# New data values for 2022.
$ws.Range("Q4").Value = 2022
$ws.Range("Q5").Value = 8725
$ws.Range("Q7").Value = 8347
$ws.Range("Q8").Value = 378

# Move the active selection to Q3, matching the saved view state.
$ws.Range("Q3").Select()
